{"js": "// Update the Dragon Spark review: new title/meta copy and refreshed\n// \"What we like\" / \"What we don't like\" bullet points.\n\nasync function replaceAll(body, findText, replaceText) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items.forEach((r) => r.insertText(replaceText, \"Replace\"));\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Title (appears both as the H1 heading and as the bold line near the end).\nawait replaceAll(\n  body,\n  \"Play Dragon Spark Free Slot | High Variance Game\",\n  \"Play Dragon Spark Free: Exciting Slot Game with Expanding Ways to Win\"\n);\n\n// \"What we like\" bullets.\nawait replaceAll(\n  body,\n  \"High-quality graphics and sound design\",\n  \"3125 ways to win with expanding rows\"\n);\nawait replaceAll(\n  body,\n  \"3125 ways to win, expanding to 16,807\",\n  \"High variance and RTP of 96.99%\"\n);\nawait replaceAll(\n  body,\n  \"High variance with a maximum payout of 10,000x\",\n  \"Cascading reels with increasing win multiplier\"\n);\nawait replaceAll(\n  body,\n  \"Cascading reels with increasing win multiplier up to 50x\",\n  \"Free Games feature with up to 16,807 ways to win\"\n);\n\n// \"What we don't like\" bullets.\nawait replaceAll(\n  body,\n  \"No progressive jackpot feature\",\n  \"Only four rows, limiting the potential for bigger wins\"\n);\nawait replaceAll(\n  body,\n  \"Bonus symbols are not easy to trigger\",\n  \"High variance may not appeal to players looking for consistent small wins\"\n);\n\n// Italic meta description line at the very end.\nawait replaceAll(\n  body,\n  \"Read our review of Dragon Spark, a high variance slot game with 3125 ways to win, cascading reels, and Free Games feature. Play for free today!\",\n  \"Play Dragon Spark free: Experience thrilling gameplay with expanding rows and multiple bonus features.\"\n);\n", "ps1": "# Update the Dragon Spark review: new title/meta copy and refreshed\n# \"What we like\" / \"What we don't like\" bullet points.\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# Title (appears both as the H1 heading and as the bold line near the end).\nReplace-AllText \"Play Dragon Spark Free Slot | High Variance Game\" \"Play Dragon Spark Free: Exciting Slot Game with Expanding Ways to Win\"\n\n# \"What we like\" bullets.\nReplace-AllText \"High-quality graphics and sound design\" \"3125 ways to win with expanding rows\"\nReplace-AllText \"3125 ways to win, expanding to 16,807\" \"High variance and RTP of 96.99%\"\nReplace-AllText \"High variance with a maximum payout of 10,000x\" \"Cascading reels with increasing win multiplier\"\nReplace-AllText \"Cascading reels with increasing win multiplier up to 50x\" \"Free Games feature with up to 16,807 ways to win\"\n\n# \"What we don't like\" bullets.\nReplace-AllText \"No progressive jackpot feature\" \"Only four rows, limiting the potential for bigger wins\"\nReplace-AllText \"Bonus symbols are not easy to trigger\" \"High variance may not appeal to players looking for consistent small wins\"\n\n# Italic meta description line at the very end.\nReplace-AllText \"Read our review of Dragon Spark, a high variance slot game with 3125 ways to win, cascading reels, and Free Games feature. Play for free today!\" \"Play Dragon Spark free: Experience thrilling gameplay with expanding rows and multiple bonus features.\"\n"}
